# farmer_buck.xlsx regen: replace the old "Strike#" derived K column values
# with the newly-computed K (std/mean based "s_vals") values, row by row,
# in column G of the save-data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number (as it appears in the worksheet, header is row 1) -> new K value
$sVals = @{
    2  = 1;
    3  = 1;
    4  = 1;
    5  = 2;
    6  = 0;
    7  = 3;
    8  = 1;
    9  = 1;
    10 = 1;
    11 = 2;
    12 = 0;
    13 = 0;
    14 = 3;
    15 = 0;
    16 = 2;
    17 = 1;
    18 = 2;
    19 = 1;
    20 = 1;
    21 = 0;
    22 = 1;
    23 = 2;
    24 = 1;
    25 = 2;
    26 = 1;
    27 = 1;
    28 = 1;
    29 = 1;
    30 = 1;
    31 = 1;
    32 = 1;
    33 = 1;
    34 = 0;
    35 = 2;
    36 = 0;
    37 = 0;
    38 = 1;
    39 = 1;
    40 = 2;
    41 = 0;
    42 = 2;
    43 = 0;
    44 = 1;
    45 = 0;
    46 = 1;
    47 = 0;
    48 = 0;
    49 = 1;
    50 = 0;
    51 = 0;
    52 = 3;
    53 = 2;
    54 = 1;
    55 = 2;
    56 = 2;
    57 = 2;
    58 = 2;
    59 = 0;
    60 = 0;
    61 = 2;
    62 = 1;
    63 = 0;
    64 = 1;
    65 = 0;
    66 = 0;
    67 = 2;
    68 = 2;
    69 = 0;
    70 = 0;
    71 = 1;
    72 = 1;
    73 = 0;
    74 = 2;
    75 = 3;
    76 = 3;
    78 = 1;
    79 = 1;
    81 = 2;
    82 = 1
}

$kCol = 7  # column G ("K")

foreach ($row in $sVals.Keys) {
    $ws.Cells.Item($row, $kCol).Value = $sVals[$row]
}
